# Sprint 1 Assigned Tasks
# Fill in the "Assigned to" column (H) on the "Sprint 1" sheet with the
# task owners, turning on word-wrap for the longer entries, then fix up
# the "Total ideal Hours" sum so it covers the two extra rows that were
# added to the sprint, and leave the sheet scrolled/selected near the
# area that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# Write the assignments in the same order they were typed in originally so
# that duplicate text reuses the same shared entry. Turn WrapText on for
# the rows whose assignee text needs to wrap inside the column.
$ws.Range("H17:H26").Value = "Assigned to: Shadi Makdissi and Ogo-Oluwa Jesutomi Olasubulumi"
$ws.Range("H17:H26").WrapText = $true

$ws.Range("H7:H8").Value = "Assigned to: Irfan Ahmed"

$ws.Range("H14:H16").Value = "Assigned to: Alec Kurkdjian and Cong-Vinh Vu"
$ws.Range("H14:H16").WrapText = $true

$ws.Range("H10:H11").Value = "Assigned to: Alec Kurkdjian"
$ws.Range("H10:H11").WrapText = $true

$ws.Range("H9").Value = "Assigned to: Cong-Vinh Vu"

$ws.Range("H12:H13").Value = "Assigned to: Cong-Vinh Vu"
$ws.Range("H12:H13").WrapText = $true

# Two more stories (rows 25 & 26) are now part of the sprint total.
$ws.Range("C5").Formula = "=SUM(E7:E26)"

# Leave the view near the area that was just edited.
$ws.Activate() | Out-Null
$ws.Range("I7").Select() | Out-Null
